# Apply renames:
#  - "Data" sheet: every occurrence of "KV1" -> "V1" and "KV2" -> "V2"
#    (these labels repeat once per question block down column C / G)
#  - "Tabellen Open" sheet header row (also updates Table1 column names):
#       A1: "V1:    1. Dit is een normale SR Vraag"   -> "V1:    . Dit is een normale SR Vraag"
#       B1: "V2:    2. Dit is een SR Vraag + Anders, namelijk" -> "V2:    . Dit is een SR Vraag + Anders, namelijk"
#       C1: "OPEN2_5:    2. Dit is een SR Vraag + Anders, namelijk A:5" -> "OPEN2_5:    . Dit is een SR Vraag + Anders, namelijk A:5"
#       D1: "OPEN7_5:    Anders, namelijk:" -> "OPEN7_5:    nders, namelijk:"
#       E1: "OPEN8_5:    Anders, namelijk:" -> "OPEN8_5:    nders, namelijk:"
#       F1: "V18_A1:    18-1. Antwoord 1" -> "V18_A1:    8-1. Antwoord 1"
#       G1: "V18_A2:    18-2. Antwoord 2" -> "V18_A2:    8-2. Antwoord 2"
#       H1: "V18_A3:    18-3. Antwoord 3" -> "V18_A3:    8-3. Antwoord 3"
#       I1: "V19:    19. Dit is een open vraag" -> "V19:    9. Dit is een open vraag"

$wb = $excel.ActiveWorkbook

# "KV1"/"KV2" are repeated (shared string) labels used throughout the Data
# sheet's header rows (not only in row 1), so use Replace on the whole used
# range to update every occurrence at once (mirrors editing the shared
# string table entry itself).
$wsData = $wb.Worksheets.Item("Data")
$usedData = $wsData.UsedRange
$usedData.Replace("KV1", "V1", 1) | Out-Null
$usedData.Replace("KV2", "V2", 1) | Out-Null

$wsTab = $wb.Worksheets.Item("Tabellen Open")
$wsTab.Range("A1").Value = "V1:    . Dit is een normale SR Vraag"
$wsTab.Range("B1").Value = "V2:    . Dit is een SR Vraag + Anders, namelijk"
$wsTab.Range("C1").Value = "OPEN2_5:    . Dit is een SR Vraag + Anders, namelijk A:5"
$wsTab.Range("D1").Value = "OPEN7_5:    nders, namelijk:"
$wsTab.Range("E1").Value = "OPEN8_5:    nders, namelijk:"
$wsTab.Range("F1").Value = "V18_A1:    8-1. Antwoord 1"
$wsTab.Range("G1").Value = "V18_A2:    8-2. Antwoord 2"
$wsTab.Range("H1").Value = "V18_A3:    8-3. Antwoord 3"
$wsTab.Range("I1").Value = "V19:    9. Dit is een open vraag"
